$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at F:G (existing F..I shift right to H..K)
$ws.Columns("F:G").Insert()

# New header cells in the freshly inserted columns
$ws.Range("F7").Value = "Group"
$ws.Range("G7").Value = "Subgroup"

# Widen the new columns
$ws.Columns("F").ColumnWidth = 17.666666666666668
$ws.Columns("G").ColumnWidth = 27.666666666666668

# Move the active selection
[void]$ws.Range("G8").Select()
